$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.547.20'
$ws.Range("E2").Value = '  -3.85%  '

$ws.Range("D3").Value = '2.510.20'
$ws.Range("E3").Value = '  -5.18%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '579.13'
$ws.Range("E5").Value = '  -2.16%  '

$ws.Range("D6").Value = '167.09'
$ws.Range("E6").Value = '  -4.63%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("D9").Value = '2.508.96'
$ws.Range("E9").Value = '  -5.24%  '

$ws.Range("D10").Value = '0.160'
$ws.Range("E10").Value = '  -6.59%  '

$ws.Range("E12").Value = '  -4.45%  '

$ws.Range("E13").Value = '  -2.09%  '

$ws.Range("D14").Value = '2.964.56'
$ws.Range("E14").Value = '  -5.38%  '

$ws.Range("D15").Value = '69.361.35'
$ws.Range("E15").Value = '  -3.93%  '

$ws.Range("D16").Value = '0.0000176'
$ws.Range("E16").Value = '  -5.51%  '

$ws.Range("D17").Value = '24.95'
$ws.Range("E17").Value = '  -4.25%  '

$ws.Range("D18").Value = '2.505.37'
$ws.Range("E18").Value = '  -7.71%  '

$ws.Range("D19").Value = '11.49'
$ws.Range("E19").Value = '  -6.87%  '

$ws.Range("D20").Value = '7.80'
$ws.Range("E20").Value = '  -2.57%  '

$ws.Range("D21").Value = '351.51'
$ws.Range("E21").Value = '  -5.16%  '

$ws.Range("D22").Value = '3.97'
$ws.Range("E22").Value = '  -4.62%  '

$ws.Range("D23").Value = '1.99'
$ws.Range("E23").Value = '  -3.82%  '

$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").Value = '69.26'
$ws.Range("E25").Value = '  -3.11%  '

$ws.Range("D26").Value = '4.03'
$ws.Range("E26").Value = '  -5.52%  '

$ws.Range("D27").Value = '9.06'
$ws.Range("E27").Value = '  -6.39%  '

$ws.Range("D28").Value = '2.637.48'
$ws.Range("E28").Value = '  -5.31%  '

$ws.Range("E29").Value = '  +1.19%  '

$ws.Range("D30").Value = '0.0₃0906'
$ws.Range("E30").Value = '  -5.01%  '

$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  -2.15%  '

$ws.Range("D32").Value = '479.49'
$ws.Range("E32").Value = '  -3.57%  '

$ws.Range("E33").Value = '  +0.57%  '

$ws.Range("E34").Value = '  -2.82%  '

$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("E36").Value = '  -1.27%  '

$ws.Range("D37").Value = '152.69'
$ws.Range("E37").Value = '  -5.70%  '

$ws.Range("D38").Value = '18.89'
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("D39").Value = '18.59'
$ws.Range("E39").Value = '  -4.01%  '

$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("E41").Value = '  -2.83%  '

$ws.Range("E42").Value = '  -2.93%  '

$ws.Range("E43").Value = '  -5.95%  '

$ws.Range("D44").Value = '1.17'
$ws.Range("E44").Value = '  -13.40%  '

$ws.Range("E45").Value = '  -8.41%  '

$ws.Range("D46").Value = '38.15'
$ws.Range("E46").Value = '  -2.63%  '

$ws.Range("D47").Value = '144.33'
$ws.Range("E47").Value = '  -6.22%  '

$ws.Range("E48").Value = '  -3.43%  '

$ws.Range("D49").Value = '0.532'
$ws.Range("E49").Value = '  -3.56%  '

$ws.Range("E50").Value = '  -5.07%  '

$ws.Range("E51").Value = '  -2.43%  '
